$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26 corresponds to "Open year" 2024.
# Update Energy Storage (column C) and Solar (column E) capacity values.
$ws.Range("C26").Value = 23.04
$ws.Range("E26").Value = 86.33199999999999
